# Increments a hex string (as text) by a given non-negative integer amount,
# preserving its original length (handles carries across all digits).
# NOTE: loop variable here is $j (not $i) to avoid colliding with the
# caller's own $i loop variable, since this runtime does not give
# functions their own variable scope for loop counters.
function Increment-Hex([string]$hex, [int]$amount) {
    $digits = $hex.ToCharArray()
    $carry = $amount
    for ($j = $digits.Length - 1; $j -ge 0 -and $carry -gt 0; $j--) {
        $d = [Convert]::ToInt32([string]$digits[$j], 16)
        $d = $d + $carry
        $carry = [Math]::Floor($d / 16)
        $d = $d % 16
        $digits[$j] = "0123456789ABCDEF"[$d]
    }
    return [string]::Join("", $digits)
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$baseEpcHex = "30300E890A01870077359401"
$baseSerial = 2000000001

for ($i = 0; $i -lt 150; $i++) {
    $row = $i + 2
    $epc = Increment-Hex $baseEpcHex $i
    $serial = $baseSerial + $i
    $ws.Cells.Item($row, 2).Value = $epc
    $ws.Cells.Item($row, 4).Value = "'" + [string]$serial
}
